# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12905
$ws1.Range("F6").Value  = 75
$ws1.Range("F7").Value  = 44
$ws1.Range("F9").Value  = 15
$ws1.Range("F10").Value = 12848
$ws1.Range("F11").Value = 282
$ws1.Range("F12").Value = 38
$ws1.Range("F13").Value = 8688
$ws1.Range("F14").Value = 7684
$ws1.Range("F15").Value = 195
$ws1.Range("F16").Value = 111
$ws1.Range("F19").Value = 983
$ws1.Range("F20").Value = 9
$ws1.Range("F24").Value = 18

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12905
$ws4.Range("F7").Value  = 75
$ws4.Range("F8").Value  = 44
$ws4.Range("F10").Value = 15
$ws4.Range("F11").Value = 12848
$ws4.Range("F12").Value = 282
$ws4.Range("F13").Value = 38
$ws4.Range("F14").Value = 8688
$ws4.Range("F15").Value = 7684
$ws4.Range("F16").Value = 195
$ws4.Range("F17").Value = 111
$ws4.Range("F20").Value = 983
$ws4.Range("F21").Value = 9
$ws4.Range("F26").Value = 18
